# Updated cryptos list on Wed May 29 12:47:03 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for every
# coin row, and swaps the Hedera / dogwifhat rows (37 <-> 38) back into their
# new ranking order (with their own refreshed price/volume too).
#
# Values such as "0.460", "36.80", "1.00" and "8.70" look like plain decimal
# numbers, so a bare .Value assignment would let Excel coerce them to the
# Number type and silently drop the significant trailing zero (0.460 -> 0.46).
# Prefixing those few with a leading apostrophe forces them to be stored as
# text, exactly like a user typing '0.460 into the cell in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = '67.746.68'
$ws.Range("E2").Value  = '  -0.97%  '

$ws.Range("D3").Value  = '3.804.36'
$ws.Range("E3").Value  = '  -2.53%  '

$ws.Range("E4").Value  = '  -0.16%  '

$ws.Range("D5").Value  = '598.26'
$ws.Range("E5").Value  = '  -0.66%  '

$ws.Range("D6").Value  = '167.77'
$ws.Range("E6").Value  = '  -2.02%  '

$ws.Range("D7").Value  = '3.801.64'
$ws.Range("E7").Value  = '  -2.57%  '

$ws.Range("E8").Value  = '  +0.13%  '

$ws.Range("D9").Value  = '0.529'
$ws.Range("E9").Value  = '  -0.35%  '

$ws.Range("D10").Value = '0.164'
$ws.Range("E10").Value = '  -0.80%  '

$ws.Range("D11").Value = '6.48'
$ws.Range("E11").Value = '  +0.57%  '

$ws.Range("D12").Value = "'0.460"
$ws.Range("E12").Value = '  -0.03%  '

$ws.Range("D13").Value = '0.0000268'
$ws.Range("E13").Value = '  +2.89%  '

$ws.Range("D14").Value = "'36.80"
$ws.Range("E14").Value = '  -1.38%  '

$ws.Range("D15").Value = '4.442.02'
$ws.Range("E15").Value = '  -2.64%  '

$ws.Range("D16").Value = '3.840.00'
$ws.Range("E16").Value = '  -1.52%  '

$ws.Range("D17").Value = '18.83'
$ws.Range("E17").Value = '  +3.29%  '

$ws.Range("D18").Value = '67.762.66'
$ws.Range("E18").Value = '  -1.11%  '

$ws.Range("D19").Value = '7.33'
$ws.Range("E19").Value = '  -0.63%  '

$ws.Range("E20").Value = '  +0.06%  '

$ws.Range("D21").Value = '10.66'
$ws.Range("E21").Value = '  -1.74%  '

$ws.Range("D22").Value = '466.52'
$ws.Range("E22").Value = '  -1.17%  '

$ws.Range("D23").Value = '0.737'
$ws.Range("E23").Value = '  -0.38%  '

$ws.Range("D24").Value = '0.0000149'
$ws.Range("E24").Value = '  -9.02%  '

$ws.Range("D25").Value = '83.56'
$ws.Range("E25").Value = '  -0.12%  '

$ws.Range("D26").Value = '2.25'
$ws.Range("E26").Value = '  +0.21%  '

$ws.Range("D27").Value = '12.17'
$ws.Range("E27").Value = '  +0.10%  '

$ws.Range("D28").Value = '10.27'
$ws.Range("E28").Value = '  +3.03%  '

$ws.Range("E29").Value = '  -0.11%  '

$ws.Range("D30").Value = '2.91'
$ws.Range("E30").Value = '  -2.01%  '

$ws.Range("D31").Value = '3.956.96'
$ws.Range("E31").Value = '  -2.52%  '

$ws.Range("D32").Value = '7.66'
$ws.Range("E32").Value = '  -2.81%  '

$ws.Range("D33").Value = '2.26'
$ws.Range("E33").Value = '  -2.45%  '

$ws.Range("D34").Value = '30.54'
$ws.Range("E34").Value = '  -2.72%  '

$ws.Range("D35").Value = '9.25'
$ws.Range("E35").Value = '  -2.08%  '

$ws.Range("D36").Value = '3.768.56'
$ws.Range("E36").Value = '  -2.78%  '

# Rows 37/38 swap: Hedera moves up to rank 35 (row 37), dogwifhat drops to
# rank 36 (row 38) - each also gets its own refreshed price/volume.
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '0.105'
$ws.Range("E37").Value = '  -0.09%  '

$ws.Range("B38").Value = 'dogwifhat'
$ws.Range("C38").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D38").Value = '3.76'
$ws.Range("E38").Value = '  +1.88%  '

$ws.Range("D39").Value = '5.93'
$ws.Range("E39").Value = '  +0.10%  '

$ws.Range("D40").Value = '0.138'
$ws.Range("E40").Value = '  -1.77%  '

$ws.Range("E41").Value = '  -2.61%  '

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = '  -0.14%  '

$ws.Range("D43").Value = '0.316'
$ws.Range("E43").Value = '  +0.87%  '

# Row 44 (USDe) is unchanged in this update.

$ws.Range("D45").Value = '1.96'
$ws.Range("E45").Value = '  -1.39%  '

$ws.Range("D46").Value = "'8.70"
$ws.Range("E46").Value = '  +0.92%  '

$ws.Range("D47").Value = '46.32'
$ws.Range("E47").Value = '  -1.88%  '

$ws.Range("D48").Value = '403.97'
$ws.Range("E48").Value = '  -5.26%  '

$ws.Range("D49").Value = '0.000277'
$ws.Range("E49").Value = '  -10.14%  '

$ws.Range("D50").Value = '142.43'
$ws.Range("E50").Value = '  -0.55%  '

$ws.Range("D51").Value = '0.0357'
$ws.Range("E51").Value = '  -0.33%  '
